$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update values in row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1

# Update values in row 19
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0

# Move/activate the sheet and set the selection to E18
$ws.Activate()
$ws.Range("E18").Select()
